# add a "users" column to the "project hours" sheet, listing the users
# that worked on each project.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

$ws.Range("E1").Value = "users"

$users = @(
    "['Daniel Olivas Hernandez', 'Berk Cagilci', 'Yhoas Olivas Hernandez']",
    "['Jiyang Chen']",
    "['Jonathan Hoff']",
    "['Alex Hill']",
    "['Chenghao DUAN', 'Arun Lakshmanan']",
    "['Zhenghe Shangguan']"
)

for ($i = 0; $i -lt $users.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $users[$i]
}

# match the header formatting already used by the other headers (bold,
# bordered, centered) on row 1.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
